$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the four "WD*" style text timestamps (shifted from September to October)
$ws.Range("R2").Value = "2017-10-10 10:00:00"
$ws.Range("S2").Value = "2017-10-10 16:00:00"
$ws.Range("L2").Value = "2017-10-21 10:00:00"
$ws.Range("M2").Value = "2017-10-21 13:00:00"

# Update the real date cell (PreferredStartDate-ish), which moves from 2017-09-21 to 2017-10-21
$ws.Range("Q2").Value = (Get-Date -Year 2017 -Month 10 -Day 21 -Hour 0 -Minute 0 -Second 0)

# Update the view: scroll so column I is at the top-left, and select M2
$ws.Range("M2").Select()
$excel.ActiveWindow.ScrollColumn = 9
